# Update cryptocurrency price/volume data in the worksheet.
# All D (Price) and E (Volume 1h) column values are stored as literal
# text strings (not numbers), so we force each target cell to Text
# format before assigning, then reset the style back to Normal so no
# stray cell-style/number-format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "20.551.11"
Set-TextValue "D3" "1.469.06"
Set-TextValue "E3" "  +2.15%  "
Set-TextValue "E4" "  +0.31%  "
Set-TextValue "D5" "0.9577"
Set-TextValue "E5" "  +4.82%  "
Set-TextValue "D6" "276.71"
Set-TextValue "E6" "  +0.61%  "
Set-TextValue "D7" "0.3561"
Set-TextValue "E7" "  -1.75%  "
Set-TextValue "D8" "0.3064"
Set-TextValue "E8" "  -0.33%  "
Set-TextValue "D9" "1.085"
Set-TextValue "E9" "  +6.70%  "
Set-TextValue "D10" "39.46"
Set-TextValue "E10" "  +2.19%  "
Set-TextValue "D11" "0.06630"
Set-TextValue "E11" "  +2.24%  "
Set-TextValue "E12" "  +0.35%  "
Set-TextValue "D13" "5.462"
Set-TextValue "E13" "  +2.35%  "
Set-TextValue "D14" "18.07"
Set-TextValue "E14" "  +3.56%  "
Set-TextValue "D15" "6.159"
Set-TextValue "E15" "  +2.18%  "
Set-TextValue "D16" "0.9588"
Set-TextValue "E16" "  +3.07%  "
Set-TextValue "D17" "0.00001020"
Set-TextValue "E17" "  +1.07%  "
Set-TextValue "D18" "1.468.25"
Set-TextValue "E18" "  +1.98%  "
Set-TextValue "D19" "0.05961"
Set-TextValue "E19" "  +5.99%  "
Set-TextValue "D20" "68.86"
Set-TextValue "E20" "  +2.21%  "
Set-TextValue "D21" "5.477"
Set-TextValue "E21" "  +1.60%  "
Set-TextValue "D22" "14.52"
Set-TextValue "E22" "  +2.16%  "
Set-TextValue "D23" "11.28"
Set-TextValue "E23" "  +4.55%  "
Set-TextValue "E24" "  +0.75%  "
Set-TextValue "D25" "20.546.31"
Set-TextValue "E25" "  +1.46%  "
Set-TextValue "D26" "145.03"
Set-TextValue "E26" "  +5.49%  "
Set-TextValue "D27" "2.089"
Set-TextValue "E27" "  -0.23%  "
Set-TextValue "D28" "17.11"
Set-TextValue "E28" "  +1.43%  "
Set-TextValue "D29" "1.630.74"
Set-TextValue "E29" "  +2.41%  "
Set-TextValue "D30" "113.94"
Set-TextValue "E30" "  +3.49%  "
Set-TextValue "D31" "3.856"
Set-TextValue "E31" "  -1.87%  "
Set-TextValue "B32" "Filecoin"
Set-TextValue "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "4.916"
Set-TextValue "E32" "  +1.79%  "
Set-TextValue "B33" "Stellar"
Set-TextValue "C33" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D33" "0.07938"
Set-TextValue "E33" "  +3.94%  "
Set-TextValue "D34" "0.7967"
Set-TextValue "E34" "  -0.78%  "
Set-TextValue "D35" "1.250"
Set-TextValue "E35" "  +11.00%  "
Set-TextValue "D36" "1.453"
Set-TextValue "E36" "  -0.86%  "
Set-TextValue "D37" "0.05757"
Set-TextValue "E37" "  -0.65%  "
Set-TextValue "D38" "4.703"
Set-TextValue "E38" "  +0.86%  "
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.02028"
Set-TextValue "E39" "  +2.38%  "
Set-TextValue "B40" "Frax"
Set-TextValue "C40" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D40" "0.9590"
Set-TextValue "E40" "  +4.00%  "
Set-TextValue "D41" "10.33"
Set-TextValue "E41" "  +1.84%  "
Set-TextValue "E42" "  +0.90%  "
Set-TextValue "D43" "7.280"
Set-TextValue "E43" "  +3.22%  "
Set-TextValue "D44" "0.5249"
Set-TextValue "E44" "  +0.96%  "
Set-TextValue "E45" "  +0.90%  "
Set-TextValue "D46" "12.04"
Set-TextValue "E46" "  +2.28%  "
Set-TextValue "D47" "118.69"
Set-TextValue "E47" "  +2.03%  "
Set-TextValue "D48" "0.5178"
Set-TextValue "E48" "  +2.07%  "
Set-TextValue "D49" "1.801"
Set-TextValue "E49" "  +4.44%  "
Set-TextValue "D50" "0.06434"
Set-TextValue "E50" "  +0.54%  "
Set-TextValue "D51" "0.9930"
